$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the three new data rows (B/C/D) for rows that previously only had
# column A populated, and remove the two trailing rows ("Crear Rendicion",
# "Agregar Inventario") that are no longer part of the table.

$ws.Range("B3").Value = "OK"
$ws.Range("D3").Value = "Pag 59"
$ws.Range("C3").Value = "IG001"

$ws.Range("B4").Value = "OK"
$ws.Range("D4").Value = "Pag 81"
$ws.Range("C4").Value = "AA001"

$ws.Range("B5").Value = "OK"
$ws.Range("D5").Value = "Pag 117"
$ws.Range("C5").Value = "AA008"

# Remove old rows 6 and 7 ("Crear Rendicion" / "Agregar Inventario") entirely.
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(6).Delete()

# Update the selection to match the recorded final state.
$ws.Range("C11").Select()
